$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$full = $p.Range
$r = $d.Range($full.Start, $full.Start + 1)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>הגבלתם</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>את</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>הגודל</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>של</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>ה</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/></w:rPr><w:t>payload</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve">? </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>כדי</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>למנוע</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>מצב</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>של</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>עומס</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>או</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/></w:rPr><w:t>DOS</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>? (</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>אין</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>דרישה</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>כזאת</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>במטלה</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="eastAsia"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>אם</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="eastAsia"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>אני</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="eastAsia"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>לא</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="eastAsia"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>טועה</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>לא</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/></w:rPr><w:t>dos</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>בא</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>בחינם</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>למסמך</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>חולשות</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:u w:val="single"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:b/><w:bCs/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>חולשה טובה לפרוטטוקול</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo" w:hint="cs"/><w:b/><w:bCs/><w:u w:val="single"/><w:rtl/></w:rPr><w:t>.!!</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Heebo" w:hAnsi="Heebo" w:cs="Heebo"/><w:u w:val="single"/><w:rtl/></w:rPr><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($frag)
"done"
